# Update "想去人数" (want-to-go count) values in the data sheets.
# Both "展览" and "全部类型" worksheets carry the same data table and
# need the same updates to column F (rows 2-5).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 8740
    $ws.Range("F3").Value = 196
    $ws.Range("F4").Value = 413
    $ws.Range("F5").Value = 94
}
